# Append the new daily snapshot row (2026-02-23) to the profit data sheet,
# matching the row pattern already present in the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 91

# Column A holds the date as literal text (e.g. "02/23/2026"), not a date
# serial. Temporarily force a text number format so Excel doesn't
# auto-convert the "mm/dd/yyyy"-looking string into a date value, then
# drop back to the default formatting so the new cell matches its
# unstyled neighbours.
$dateCell = $ws.Range("A$row")
$dateCell.NumberFormat = "@"
$dateCell.Value = "02/23/2026"
$dateCell.ClearFormats()

$ws.Range("B$row").Value = 9217.540000000001
$ws.Range("C$row").Value = 0.2464356284593587
$ws.Range("D$row").Value = 0.7535643715406413
$ws.Range("E$row").Value = -337.53
$ws.Range("F$row").Value = -36.1
$ws.Range("G$row").Value = -24045.75
$ws.Range("H$row").Value = -77.59
$ws.Range("I$row").Value = -1181.32
$ws.Range("J$row").Value = -34.21
$ws.Range("K$row").Value = -25227.07
$ws.Range("L$row").Value = -73.23999999999999
